$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp update (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 15:52"

# Row 4 - Estados Unidos: updated totals
$ws.Range("B4").Value = 1097080
$ws.Range("C4").Value = 2057
$ws.Range("D4").Value = 155829
$ws.Range("E4").Value = 877338
$ws.Range("G4").Value = 57
$ws.Range("H4").Value = 63913

# Row 23 - Arabia Saudita: Casos criticos updated
$ws.Range("F23").Value = 117

# Row 42 - Serbia: Casos criticos updated
$ws.Range("F42").Value = 65

# Rows 74/75 - Azerbaiyan overtakes Camerun in total cases, so the sorted
# table swaps their order. Row 74 now holds Azerbaiyan's refreshed figures
# and row 75 holds Camerun's (previously row-74) figures.
$ws.Range("A74").Value = "Azerbaiyan"
$ws.Range("B74").Value = 1854
$ws.Range("C74").Value = 50
$ws.Range("D74").Value = 1365
$ws.Range("E74").Value = 464
$ws.Range("F74").Value = 17
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 25

$ws.Range("A75").Value = "Camerun"
$ws.Range("B75").Value = 1832
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 934
$ws.Range("E75").Value = 837
$ws.Range("F75").Value = 12
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 61

# Row 110 - Georgia: Casos activos / Recuperados updated
$ws.Range("D110").Value = 207
$ws.Range("E110").Value = 353
